# Update attendance/price figures to match the latest scraped data snapshot
# (commit: "Update gh-pages to output generated at 456a3b4")
$wb = $excel.ActiveWorkbook

$ws1 = $wb.Worksheets.Item("展览")
$ws2 = $wb.Worksheets.Item("演出")
$ws3 = $wb.Worksheets.Item("本地生活")
$ws4 = $wb.Worksheets.Item("全部类型")

# --- 展览 ---
$ws1.Range("F2").Value = 1527
$ws1.Range("F5").Value = 7892
$ws1.Range("F6").Value = 4906
$ws1.Range("G6").Value = 80
$ws1.Range("F7").Value = 7197
$ws1.Range("G7").Value = 75
$ws1.Range("F8").Value = 299
$ws1.Range("G8").Value = 70
$ws1.Range("F9").Value = 1534
$ws1.Range("G9").Value = 70
$ws1.Range("F10").Value = 883
$ws1.Range("F12").Value = 74
$ws1.Range("F15").Value = 569
$ws1.Range("F16").Value = 33
$ws1.Range("F17").Value = 246
$ws1.Range("F20").Value = 1247
$ws1.Range("F24").Value = 1272
$ws1.Range("F25").Value = 55
$ws1.Range("F30").Value = 228
$ws1.Range("F31").Value = 1026
$ws1.Range("F33").Value = 16
$ws1.Range("F34").Value = 155
$ws1.Range("F35").Value = 138
$ws1.Range("F42").Value = 111
$ws1.Range("F44").Value = 1208
$ws1.Range("F45").Value = 612

# --- 演出 ---
$ws2.Range("G6").Value = "不可售"
$ws2.Range("F10").Value = 139
$ws2.Range("G10").Value = "不可售"
$ws2.Range("F25").Value = 643
$ws2.Range("F29").Value = 39
$ws2.Range("F32").Value = 880
$ws2.Range("F34").Value = 1006
$ws2.Range("F42").Value = 148

# --- 本地生活 ---
$ws3.Range("F6").Value = 704
$ws3.Range("F7").Value = 202
$ws3.Range("F8").Value = 109
$ws3.Range("F9").Value = 1774
$ws3.Range("F10").Value = 2674

# --- 全部类型 ---
$ws4.Range("F3").Value = 1527
$ws4.Range("F6").Value = 704
$ws4.Range("F7").Value = 7892
$ws4.Range("F8").Value = 202
$ws4.Range("F9").Value = 4906
$ws4.Range("G9").Value = 80
$ws4.Range("F10").Value = 7197
$ws4.Range("G10").Value = 75
$ws4.Range("F11").Value = 299
$ws4.Range("G11").Value = 70
$ws4.Range("F12").Value = 1534
$ws4.Range("G12").Value = 70
$ws4.Range("F13").Value = 883
$ws4.Range("F14").Value = 109
$ws4.Range("F16").Value = 1774
$ws4.Range("F17").Value = 2674
$ws4.Range("F22").Value = 33
$ws4.Range("F23").Value = 246
$ws4.Range("F24").Value = 1247
$ws4.Range("F25").Value = 643
$ws4.Range("F27").Value = 1272
$ws4.Range("F30").Value = 228
$ws4.Range("F33").Value = 39
$ws4.Range("F34").Value = 880
$ws4.Range("F35").Value = 155
$ws4.Range("F37").Value = 138
$ws4.Range("F38").Value = 1006
$ws4.Range("F42").Value = 111
$ws4.Range("F45").Value = 612
